{"js": "// Add spacing-after=0 to the \"PROFESSIONAL EXPERIENCE\" heading paragraph,\n// insert a new empty paragraph styled \"Travis Normal Bold\" right after it,\n// and update the \"Travis Normal Bold 2\" style definition to also carry\n// spacing-after=0 in its paragraph formatting.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\n// Locate the \"PROFESSIONAL EXPERIENCE\" heading paragraph (falls back to the\n// first paragraph of the body if the text can't be matched for some reason).\nlet heading = paragraphs.items.find(\n  (p) => p.text.trim() === \"PROFESSIONAL EXPERIENCE\"\n);\nif (!heading) {\n  heading = paragraphs.items[0];\n}\n\n// Paragraph-level override matching <w:spacing w:after=\"0\"/> in <w:pPr>.\nheading.spaceAfter = 0;\n\n// New, empty paragraph right after the heading, styled \"Travis Normal Bold\".\nconst newPara = heading.insertParagraph(\"\", \"After\");\nnewPara.style = \"Travis Normal Bold\";\n\n// The \"Travis Normal Bold 2\" style definition also gains spacing-after=0.\nconst styles = context.document.getStyles();\nconst travisBold2 = styles.getByNameOrNullObject(\"Travis Normal Bold 2\");\nawait context.sync();\n\nif (!travisBold2.isNullObject) {\n  travisBold2.paragraphFormat.spaceAfter = 0;\n}\n\nawait context.sync();\n", "ps1": "# Add spacing-after=0 to the \"PROFESSIONAL EXPERIENCE\" heading paragraph,\n# insert a new empty paragraph styled \"Travis Normal Bold\" right after it,\n# and update the \"Travis Normal Bold 2\" style definition to also carry\n# spacing-after=0 in its paragraph formatting.\n\n$d = $word.ActiveDocument\n\n# Locate the \"PROFESSIONAL EXPERIENCE\" heading paragraph (falls back to the\n# first paragraph of the document if the text can't be matched for some\n# reason).\n$heading = $null\nforeach ($p in $d.Paragraphs) {\n  if ($p.Range.Text.Trim() -eq \"PROFESSIONAL EXPERIENCE\") {\n    $heading = $p\n    break\n  }\n}\nif ($heading -eq $null) {\n  $heading = $d.Paragraphs(1)\n}\n\n# Paragraph-level override matching <w:spacing w:after=\"0\"/> in <w:pPr>.\n$heading.Range.ParagraphFormat.SpaceAfter = 0\n\n# New, empty paragraph right after the heading, styled \"Travis Normal Bold\".\n$r = $heading.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n$newPara = $heading.Next()\n$newPara.Style = \"Travis Normal Bold\"\n\n# The \"Travis Normal Bold 2\" style definition also gains spacing-after=0.\ntry {\n  $travisBold2 = $d.Styles(\"Travis Normal Bold 2\")\n  $travisBold2.ParagraphFormat.SpaceAfter = 0\n} catch {\n}\n"}
